# Apply edits to Scope_1_stationary_fuel workbook
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename "Gasoline" -> "Motor Gasoline" in B5
$ws.Range("B5").Value = "Motor Gasoline"

# Update CO2 factor for Natural Gas row (C4): 53 -> 53.06
$ws.Range("C4").Value = 53.06

# Add new "Combustion Type" column (J) with header and "Stationary" values
# Copy the formatting from an existing header cell (F3) onto J3, then set its text
$ws.Range("F3").Copy()
$ws.Range("J3").PasteSpecial(-4122)
$ws.Range("J3").Value = "Combustion Type"

$ws.Range("J4").Value = "Stationary"
$ws.Range("J5").Value = "Stationary"
$ws.Range("J6").Value = "Stationary"

# Set column J width to match other similar columns
$ws.Columns.Item(10).ColumnWidth = 12.44140625

# Update selection to the new column
$ws.Range("J3:J6").Select()
